$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $s = $r.Style
    $r.Value = "'" + $val
    $r.Style = $s
}

Set-TextValue "D2" "30.328.44"
$ws.Range("E2").Value = "  -3.50%  "
Set-TextValue "D3" "1.931.08"
$ws.Range("E3").Value = "  -3.79%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "248.77"
$ws.Range("E5").Value = "  -3.94%  "
Set-TextValue "D6" "0.7249"
$ws.Range("E6").Value = "  -6.05%  "
Set-TextValue "D7" "0.9996"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue "D8" "0.3277"
$ws.Range("E8").Value = "  -8.62%  "
Set-TextValue "D9" "27.20"
$ws.Range("E9").Value = "  -3.88%  "
Set-TextValue "D10" "0.06807"
$ws.Range("E10").Value = "  -3.62%  "
Set-TextValue "D11" "0.8035"
$ws.Range("E11").Value = "  -4.33%  "
Set-TextValue "D12" "0.08048"
$ws.Range("E12").Value = "  -0.62%  "
Set-TextValue "D13" "1.929.21"
$ws.Range("E13").Value = "  -3.84%  "
Set-TextValue "D14" "5.416"
Set-TextValue "D15" "94.78"
$ws.Range("E15").Value = "  -6.45%  "
Set-TextValue "D16" "14.49"
$ws.Range("E16").Value = "  -1.15%  "
Set-TextValue "D17" "30.308.57"
$ws.Range("E17").Value = "  -3.60%  "
Set-TextValue "D18" "254.27"
$ws.Range("E18").Value = "  -7.44%  "
Set-TextValue "D19" "0.000008010"
$ws.Range("E19").Value = "  +0.56%  "
Set-TextValue "D20" "5.826"
$ws.Range("E20").Value = "  -1.98%  "
Set-TextValue "D21" "2.186.66"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("E22").Value = "  +0.01%  "
Set-TextValue "D23" "1.000"
$ws.Range("E23").Value = "  +0.06%  "
Set-TextValue "D24" "6.866"
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("E25").Value = "  -4.57%  "
Set-TextValue "D26" "159.41"
$ws.Range("E26").Value = "  -2.86%  "
Set-TextValue "D27" "2.394"
$ws.Range("E27").Value = "  +0.61%  "
Set-TextValue "D28" "0.1341"
$ws.Range("E28").Value = "  -8.80%  "
Set-TextValue "D29" "19.06"
$ws.Range("E29").Value = "  -5.31%  "
Set-TextValue "D30" "1.557"
$ws.Range("E30").Value = "  -4.36%  "
Set-TextValue "D31" "1.337"
$ws.Range("E31").Value = "  -1.49%  "
Set-TextValue "D32" "4.397"
$ws.Range("E32").Value = "  -5.00%  "
Set-TextValue "D33" "4.186"
$ws.Range("E33").Value = "  -4.94%  "
Set-TextValue "D34" "0.05064"
$ws.Range("E34").Value = "  -2.87%  "
Set-TextValue "D35" "1.217"
$ws.Range("E35").Value = "  -1.35%  "
Set-TextValue "D36" "0.7383"
$ws.Range("E36").Value = "  -3.07%  "
Set-TextValue "D37" "2.752"
$ws.Range("E37").Value = "  -1.63%  "
Set-TextValue "D38" "0.01970"
$ws.Range("E38").Value = "  -2.48%  "
Set-TextValue "D39" "2.828"
$ws.Range("E39").Value = "  -4.40%  "
Set-TextValue "D40" "6.593"
$ws.Range("E40").Value = "  -1.93%  "
Set-TextValue "D41" "79.13"
$ws.Range("E41").Value = "  -1.37%  "
Set-TextValue "D42" "0.4455"
$ws.Range("E42").Value = "  -6.16%  "
Set-TextValue "D43" "1.994"
$ws.Range("E43").Value = "  -9.08%  "
Set-TextValue "D44" "0.9996"
$ws.Range("E44").Value = "  -0.06%  "
Set-TextValue "D45" "0.8338"
$ws.Range("E45").Value = "  -3.22%  "
Set-TextValue "D46" "101.85"
Set-TextValue "D47" "9.741"
$ws.Range("E47").Value = "  -2.30%  "
Set-TextValue "D48" "7.281"
$ws.Range("E48").Value = "  -5.04%  "
Set-TextValue "D49" "36.37"
$ws.Range("E49").Value = "  -1.93%  "
Set-TextValue "D50" "0.05940"
$ws.Range("E50").Value = "  -0.64%  "
Set-TextValue "D51" "0.4069"
$ws.Range("E51").Value = "  -6.90%  "
